$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.073.07"
$ws.Range("E2").Value = "  -3.94%  "
$ws.Range("D3").Value = "2.258.29"
$ws.Range("E3").Value = "  -4.76%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'489.61"
$ws.Range("E5").Value = "  -2.18%  "
$ws.Range("D6").Value = "'126.72"
$ws.Range("E6").Value = "  -2.59%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  -4.30%  "
$ws.Range("D9").Value = "2.259.07"
$ws.Range("E9").Value = "  -4.79%  "
$ws.Range("E10").Value = "  -6.45%  "
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("D12").Value = "'4.73"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("D14").Value = "2.657.31"
$ws.Range("E14").Value = "  -4.78%  "
$ws.Range("D15").Value = "'21.02"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("D16").Value = "54.016.78"
$ws.Range("E16").Value = "  -3.91%  "
$ws.Range("E17").Value = "  -2.70%  "
$ws.Range("D18").Value = "2.248.83"
$ws.Range("E18").Value = "  -4.39%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'3.98"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'9.61"
$ws.Range("E20").Value = "  -4.28%  "
$ws.Range("D21").Value = "'300.83"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").Value = "'6.11"
$ws.Range("E22").Value = "  -1.98%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'63.70"
$ws.Range("E24").Value = "  -2.12%  "
$ws.Range("D25").Value = "'1.01"
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").Value = "'0.366"
$ws.Range("E26").Value = "  -0.46%  "
$ws.Range("E27").Value = "  -2.99%  "
$ws.Range("D28").Value = "'7.08"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").Value = "'169.69"
$ws.Range("E29").Value = "  -1.43%  "
$ws.Range("D30").Value = "0.0₃0689"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("D31").Value = "'1.60"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'5.74"
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "'0.996"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").Value = "'1.07"
$ws.Range("E35").Value = "  -1.02%  "
$ws.Range("D36").Value = "'17.41"
$ws.Range("E36").Value = "  -0.80%  "
$ws.Range("E37").Value = "  -1.01%  "
$ws.Range("D38").Value = "'0.835"
$ws.Range("E38").Value = "  +5.56%  "
$ws.Range("D39").Value = "'3.59"
$ws.Range("E39").Value = "  -4.82%  "
$ws.Range("D40").Value = "'35.71"
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").Value = "'0.367"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").Value = "'1.37"
$ws.Range("E42").Value = "  -2.45%  "
$ws.Range("D43").Value = "'3.30"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "'122.23"
$ws.Range("E44").Value = "  -6.61%  "
$ws.Range("D45").Value = "'4.68"
$ws.Range("E45").Value = "  -1.63%  "
$ws.Range("D46").Value = "'0.0878"
$ws.Range("E46").Value = "  -2.80%  "
$ws.Range("D47").Value = "'0.538"
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("D48").Value = "'237.67"
$ws.Range("E48").Value = "  -1.22%  "
$ws.Range("D50").Value = "'0.0203"
$ws.Range("E50").Value = "  -2.68%  "
$ws.Range("D51").Value = "'16.26"
$ws.Range("E51").Value = "  -3.07%  "

# Reset number-looking text cells back to the default (unstyled) cell style
# now that the quote-prefix has locked them in as text, so we do not leave a
# stray style index behind.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
